$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-07-26 07:39:41"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("G3").Value = "2016-07-26 07:39:32"
$wsZhCn.Range("J3").Value = "2016-07-26 07:40:26"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("J3").Value = "2016-07-26 07:40:42"
